$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The sheet originally had its column-A labels offset by one row from
# the B/C data next to them (a data-entry bug). This edit inserts a
# new row 13 (so the "Docentes responsaveis:" value moves into its own
# row under the existing label) and re-populates the B/C content for
# every row so each label lines up with the text that belongs to it,
# finally adding a new row 22 for the "Bibliografia:" reference text.
# ------------------------------------------------------------------

# 1. Insert a new row at 13 - this shifts old rows 13-21 down to 14-22
#    (row heights travel with the shifted rows automatically).
$ws.Rows("13:13").Insert()

# Give the new row 13 the same B/C formatting as a normal data row
# (style only, no borders/heights), then drop the leftover A13 cell
# that Insert() cloned down from row 12.
$ws.Range("B10:C10").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)
$ws.Range("A13").Clear()

# 2. Row 10 - Objetivos: replace the (wrongly duplicated) teacher name
#    with the real objectives text.
$ws.Range("B10").Value = "Falar basicamente sobre ciência; Técnica; Tecnologia; Engenharia; Pesquisa; Descobertas e Invenções."
$ws.Range("C10").Value = "Falar basicamente sobre ciência; Técnica; Tecnologia; Engenharia; Pesquisa; Descobertas e Invenções."

# 3. Row 13 - Docentes responsáveis: value (now correctly under its own
#    label in row 12, with no label of its own).
$ws.Range("B13").Value = "6376612 - Daisy Rafaela da Silva"
$ws.Range("C13").Value = "6376612 - Daisy Rafaela da Silva"

# 4. Row 14 - Programa resumido: replace "Semestral" with the real
#    summary text.
$ws.Range("B14").Value = "Ciência ontem e hoje; Positivismo e o Neopositivismo; Física moderna e seus pensadores; As escolas de Engenharia; Preparação de monografias."
$ws.Range("C14").Value = "Ciência ontem e hoje; Positivismo e o Neopositivismo; Física moderna e seus pensadores; As escolas de Engenharia; Preparação de monografias."

# 5. Row 15 - Short syllabus: no longer carries the (wrong) date value.
$ws.Range("B15").Clear()
$ws.Range("C15").Clear()

# 6. Row 16 - Programa: gets the full syllabus text.
$ws.Range("B16").Value = "1. Ciência, técnica, tecnologia e engenharia 2. Ciência e técnica na Idade Antiga 3. Ciência e técnica na Idade Média . 4. Ciência e técnica na Idade Moderna. 5. Ciência e técnica na Idade Contemporânea 6. Metodologia Científica: Aristóteles, Galileu e Descartes. 7. Positivismo e neopositivismo, COMTE, Popper, Kuhn. 8. Definição, medidas, leis e teoria física 9. Física Moderna e realidade. 10. Teoria dos Modelos. Newton, Poincaré, Lorentz, Einstein 11. Pesquisa de causas. Leis estatísticas, determinismo e acaso 12. Ciência teórica e ciência experimental 13. Pesquisas, descobertas e invenções 14. As escolas de engenharia, formação das escolas, Escola de Engenharia de Lorena /EEL/USP. 15. Engenharia, matemática e física 16. Organização da pesquisa tecnológica 17. Preparação de Monografias tecnológicas."
$ws.Range("C16").Value = "1. Ciência, técnica, tecnologia e engenharia 2. Ciência e técnica na Idade Antiga 3. Ciência e técnica na Idade Média . 4. Ciência e técnica na Idade Moderna. 5. Ciência e técnica na Idade Contemporânea 6. Metodologia Científica: Aristóteles, Galileu e Descartes. 7. Positivismo e neopositivismo, COMTE, Popper, Kuhn. 8. Definição, medidas, leis e teoria física 9. Física Moderna e realidade. 10. Teoria dos Modelos. Newton, Poincaré, Lorentz, Einstein 11. Pesquisa de causas. Leis estatísticas, determinismo e acaso 12. Ciência teórica e ciência experimental 13. Pesquisas, descobertas e invenções 14. As escolas de engenharia, formação das escolas, Escola de Engenharia de Lorena /EEL/USP. 15. Engenharia, matemática e física 16. Organização da pesquisa tecnológica 17. Preparação de Monografias tecnológicas."

# 7. Row 17 (Syllabus:) and Row 18 (Avaliação:) no longer carry any
#    B/C content (they were only ever wrongly-shifted duplicates).
$ws.Range("B17").Clear()
$ws.Range("C17").Clear()
$ws.Range("B18").Clear()
$ws.Range("C18").Clear()

# 8. Row 19 - Método: now correctly receives the teaching-method text
#    (previously stuck one row up, i.e. what used to be row 18).
$ws.Range("B19").Value = "Aulas expositivas em nível de conferência. Estudo de casos significativos da história da ciência e da engenharia. Debate participativo em torno de questões relevantes."
$ws.Range("C19").Value = "Aulas expositivas em nível de conferência. Estudo de casos significativos da história da ciência e da engenharia. Debate participativo em torno de questões relevantes."

# 9. Row 20 - Critério: receives the grading-formula text.
$ws.Range("B20").Value = "A = (P + T)/ 2 Onde: P = média das provas T = média dos trabalhos práticos"
$ws.Range("C20").Value = "A = (P + T)/ 2 Onde: P = média das provas T = média dos trabalhos práticos"

# 10. Row 21 - Norma de recuperação: receives the recovery-exam text.
$ws.Range("B21").Value = "RECUPERAÇÃO: 1 (uma) prova."
$ws.Range("C21").Value = "RECUPERAÇÃO: 1 (uma) prova."

# 11. Row 22 (new) - Bibliografia: label plus the reference list. Copy
#     the A/B/C formatting down from row 21 (same label/text column
#     styles used throughout the sheet) before writing the new values,
#     then restore the taller 120pt row height used by other long
#     multi-line entries (e.g. row 16/17's "Programa"/"Syllabus").
$ws.Range("A21:C21").Copy()
$ws.Range("A22:C22").PasteSpecial(-4122)
$ws.Rows("22:22").RowHeight = 120

$ws.Range("A22").Value = "Bibliografia:"
$ws.Range("B22").Value = "VARGAS, Milton. Metodologia da pesquisa tecnológica,Rio de Janeiro, Globo, 1985. SIMARD, Emile. Naturaleza y alcance del método científico, Madrim, Gredos, 1961. ROUSSEAU, Pierre. História da Ciência, Lisboa, 1963. VARGAS, Milton. Para uma filosofia da tecnologia, Ed.Alfa-Omega, São Paulo, 1994."
$ws.Range("C22").Value = "VARGAS, Milton. Metodologia da pesquisa tecnológica,Rio de Janeiro, Globo, 1985. SIMARD, Emile. Naturaleza y alcance del método científico, Madrim, Gredos, 1961. ROUSSEAU, Pierre. História da Ciência, Lisboa, 1963. VARGAS, Milton. Para uma filosofia da tecnologia, Ed.Alfa-Omega, São Paulo, 1994."

Write-Host "Edit complete"
